$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '243.17'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '8'

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '8'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.421'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '8'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05928'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '8'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.456'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '8'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.539'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '8'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8133'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '8'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9139'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '8'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1406'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '8'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07427'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '8'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03278'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '8'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03055'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '8'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09351'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '8'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.848'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '8'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001575'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '8'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04670'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '8'

$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006189'
$ws.Range("E18").Value = '17TigerCashTCH'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '8'

$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004986'
$ws.Range("E19").Value = '18HotbitTokenHTB'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '8'

$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009806'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '8'

$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0001100'
$ws.Range("E21").Value = '20NitroExNTX'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '8'

$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.605'
$ws.Range("E22").Value = '21LEOLEO'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '8'

$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.136'
$ws.Range("E23").Value = '22BTSETokenBTSE'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '8'

$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.01119'
$ws.Range("E24").Value = '23OneONEBestin24h'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '8'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3227'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '8'

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '8'

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '8'

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '8'

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '8'

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '8'

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '8'

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '8'

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '8'

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '8'

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '8'

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '8'

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '8'

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '8'

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '8'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04029'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '8'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006203'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '8'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1074'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '8'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.003001'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '8'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008720'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '8'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005243'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '8'

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '8'

$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '8'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.8162'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '8'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002262'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '8'

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '8'

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '8'
